$wb = $excel.ActiveWorkbook

# --- Data Management sheet: add new designation "Maintenance Supervisor" ---
$dataMgmt = $wb.Worksheets.Item("Data Management")
$dataMgmt.Range("E3").Value = "Maintenance Supervisor"

$matReg = $wb.Worksheets.Item("Material Registration")

# --- Update active cell selections to match the authored state ---
$matReg.Activate()
$matReg.Range("I3").Select()

$dataMgmt.Activate()
$dataMgmt.Range("H9").Select()

$matReg.Activate()
